$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append ".tif" to each filename in column A (rows 2 through 25)
for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + ".tif"
}

# Update the view: scroll so row 6 is at the top, and select A25 (last edited cell)
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select()
